$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 50
$ws.Range("C2").Value = 80
$ws.Range("B3").Value = 50
$ws.Range("C3").Value = 80
$ws.Range("B4").Value = 40
$ws.Range("C4").Value = 80
$ws.Range("B5").Value = 50
$ws.Range("C5").Value = 90
$ws.Range("B6").Value = 40
$ws.Range("C6").Value = 80
$ws.Range("B7").Value = 50
$ws.Range("C7").Value = 90

$ws.Range("C8").Select()
$excel.ActiveWindow.ScrollRow = 9
